# Archetypes.xlsx - "Minor tweaks to feats"
# Applies: four text tweaks to feat/archetype description cells, an
# AutoFilter range change (A1:K9 -> A1:M10) together with the matching
# swap of the two hidden/visible _FilterDatabase defined names, a row
# height tweak on row 3, and the resulting view/selection changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Update the 4 shared-string text cells (content tweaks to feats) ---

# G3: Auror innateDescription (Intuition)
$ws.Range("G3").Value2 = '\imp{\innateAbility} is the inherent\comma{} instinctive understanding of the minds of others possessed by an insightful and trained mind. Bypassing all \imp{Logic} and conscious reasoning\comma{} \imp{intuition} allows an \name{} to make great strides in their understanding of people and their actions by getting inside their heads and understanding the way that they think. \imp{\innateAbility} can allow an \bname{} to suddenly have a flash of insight into the motives\comma{} aims or drive of another being: If you wish to know why someone would behave in a given way\comma{} why a certain shop was robbed and not another\comma{} or where a target might head next - an \name{}''s \imp{\innateAbility} is surely the best tool'

# K3: Auror knowledgeDescription (Tracking)
$ws.Range("K3").Value2 = 'Hunting down a foe is a key part of being an \bname{}, and part of that is being able to survey a scene and see where they were, what they did, and where they''re going next. When you \imp{Track} a target you look for the trail that they have left - scuffs in the dirt, broken twigs in the forest and even more abstract trails such as an online presence or a paper trail. '

# D7: Responder experience text
$ws.Range("D7").Value2 = '\item Overcome a problem with care and  compassion before considering violence. 
\item Treat a serious wound, heal an injured ally or prevent an injury from occuring '

# G10: Warrior innateDescription (Rage)
$ws.Range("G10").Value2 = '\imp{Rage} is the deep seated anger that lies within the hearts of most people, even the most benevolent of us. A \bname{}, however, has learned to weaponise their rage, either by letting it out in an unbridled fury, or harnessing it, fuelling their cold, calculated actions. 

Whilst in combat, \imp{Rage} can be substituted for almost any physical act such as a weapon attack, and can often serve as a useful social crutch when you need to terrify someone, or need adrenaline to lift a fallen tree from an ally. '

# --- 2. Row 3 auto-height shrank slightly once the Intuition text lost a paragraph ---
$ws.Rows.Item(3).RowHeight = 189.55

# --- 3. Grow the AutoFilter range from A1:K9 to A1:M10 ---
$ws.AutoFilterMode = $false
$ws.Range("A1:M10").AutoFilter() | Out-Null

# The AutoFilter range change leaves the workbook's two "_FilterDatabase"
# defined names (one hidden, one visible leftover) pointing at the ranges
# in their original order; swap their RefersTo so the hidden one tracks
# the new/current autofilter range (A1:M10) and the stale visible one
# keeps the old range (A1:K9), matching the authored workbook.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        if (-not $n.Visible) {
            $n.RefersTo = "=Sheet1!`$A`$1:`$M`$10"
        } else {
            $n.RefersTo = "=Sheet1!`$A`$1:`$K`$9"
        }
    }
}

# --- 4. View state: freeze the first column and move the visible window
#        down so row 8 is at the top, with A8 selected in the frozen pane
#        and I9 selected in the scrolling pane. ---
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("B1").Select() | Out-Null
$win.FreezePanes = $true
$win.ScrollRow = 8
$win.ScrollColumn = 3
$ws.Range("A8").Select() | Out-Null
$ws.Range("I9").Select() | Out-Null
